$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (149) down into the two
# new rows so the new cells pick up the same styles (bold/bordered index
# column, date-formatted match-time column) as every other row.
$ws.Range("A149:V149").Copy()
$ws.Range("A150:V151").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 150 - Kyzylzhar 0 x 0 Aktobe
$ws.Range("A150").Value = 149
$ws.Range("B150").Value = "kazakhstan"
$ws.Range("C150").Value = "premier-league"
$ws.Range("D150").Value = "2023"
$ws.Range("E150").Value = 45192.45833333334
$ws.Range("F150").Value = "Kyzylzhar"
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = "Aktobe"
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 2.73
$ws.Range("K150").Value = "21/09/2023 22:13"
$ws.Range("L150").Value = 2.94
$ws.Range("M150").Value = "23/09/2023 10:57"
$ws.Range("N150").Value = 2.85
$ws.Range("O150").Value = "21/09/2023 22:13"
$ws.Range("P150").Value = 2.84
$ws.Range("Q150").Value = "23/09/2023 10:57"
$ws.Range("R150").Value = 2.33
$ws.Range("S150").Value = "21/09/2023 22:13"
$ws.Range("T150").Value = 2.43
$ws.Range("U150").Value = "23/09/2023 10:57"
$ws.Range("V150").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/kyzylzhar-aktobe/lYlcb3LR/"

# Row 151 - Maqtaaral 2 x 0 Aksu
$ws.Range("A151").Value = 150
$ws.Range("B151").Value = "kazakhstan"
$ws.Range("C151").Value = "premier-league"
$ws.Range("D151").Value = "2023"
$ws.Range("E151").Value = 45192.54166666666
$ws.Range("F151").Value = "Maqtaaral"
$ws.Range("G151").Value = 2
$ws.Range("H151").Value = "Aksu"
$ws.Range("I151").Value = 0
$ws.Range("J151").Value = 2.24
$ws.Range("K151").Value = "22/09/2023 04:12"
$ws.Range("L151").Value = 1.67
$ws.Range("M151").Value = "23/09/2023 12:59"
$ws.Range("N151").Value = 2.99
$ws.Range("O151").Value = "22/09/2023 04:12"
$ws.Range("P151").Value = 3.3
$ws.Range("Q151").Value = "23/09/2023 12:59"
$ws.Range("R151").Value = 2.73
$ws.Range("S151").Value = "22/09/2023 04:12"
$ws.Range("T151").Value = 4.19
$ws.Range("U151").Value = "23/09/2023 12:59"
$ws.Range("V151").Value = "https://www.betexplorer.com/football/kazakhstan/premier-league/maqtaaral-aksu/SI0l0PjF/"

# Column D ("temporada") holds the season as TEXT in every existing row
# ("2023"), but a plain .Value assignment of a numeric-looking string gets
# auto-coerced to a number. Route it through a formula + paste-values so
# the final cell is a genuine text value without touching NumberFormat
# (which would otherwise mint an unwanted extra cell style).
$ws.Range("D150").Formula = '="2023"'
$ws.Range("D151").Formula = '="2023"'
$ws.Range("D150:D151").Copy()
$ws.Range("D150:D151").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

Write-Host "Added rows 150 and 151"
